$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: date 44294 -> 44295
$ws.Range("D2").Value2 = 44295

# Row 3: date 44294 -> 44295; Volumen (M3) 240 -> 200
$ws.Range("D3").Value2 = 44295
$ws.Range("M3").Value2 = 200

# Row 4: date 44294 -> 44295
$ws.Range("D4").Value2 = 44295

# Row 5: date 44295 -> 44294
$ws.Range("D5").Value2 = 44294

# Row 6: date 44295 -> 44294; Volumen (M6) 200 -> 240
$ws.Range("D6").Value2 = 44294
$ws.Range("M6").Value2 = 240

# Row 7: date 44295 -> 44294
$ws.Range("D7").Value2 = 44294
